$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 26.127733
$ws.Range("H2").Value = 78.383199
$ws.Range("I2").Value = 0.2666992864894373
$ws.Range("J2").Value = 0.2666992864894374
$ws.Range("M2").Value = 54.89331066666666
$ws.Range("N2").Value = 164.679932
$ws.Range("O2").Value = 0.8971624892852424
$ws.Range("P2").Value = 0.8971624892852424
$ws.Range("Q2").Value = 1434.237764584719
$ws.Range("R2").Value = 12908.13988126247
$ws.Range("S2").Value = 0.2392725957574616
$ws.Range("T2").Value = 0.2392725957574617
$ws.Range("G3").Value = 26.127733
$ws.Range("H3").Value = 78.383199
$ws.Range("I3").Value = 0.2666992864894373
$ws.Range("J3").Value = 0.2666992864894374
$ws.Range("M3").Value = 2.309992333333333
$ws.Range("N3").Value = 6.929977
$ws.Range("O3").Value = 0.03775393480250816
$ws.Range("P3").Value = 0.03775393480250816
$ws.Range("Q3").Value = 60.35486291738034
$ws.Range("R3").Value = 543.193766256423
$ws.Range("S3").Value = 0.01006894747399766
$ws.Range("T3").Value = 0.01006894747399766
$ws.Range("G4").Value = 26.127733
$ws.Range("H4").Value = 78.383199
$ws.Range("I4").Value = 0.2666992864894373
$ws.Range("J4").Value = 0.2666992864894374
$ws.Range("M4").Value = 3.982169333333333
$ws.Range("N4").Value = 11.946508
$ws.Range("O4").Value = 0.06508357591224938
$ws.Range("P4").Value = 0.06508357591224936
$ws.Range("Q4").Value = 104.0450571021213
$ws.Range("R4").Value = 936.4055139190921
$ws.Range("S4").Value = 0.01735774325797804
$ws.Range("T4").Value = 0.01735774325797804
$ws.Range("I5").Value = 0.2440410104700376
$ws.Range("J5").Value = 0.2440410104700377
$ws.Range("M5").Value = 54.89331066666666
$ws.Range("N5").Value = 164.679932
$ws.Range("O5").Value = 0.8971624892852424
$ws.Range("P5").Value = 0.8971624892852424
$ws.Range("Q5").Value = 1312.387588023806
$ws.Range("R5").Value = 11811.48829221425
$ws.Range("S5").Value = 0.2189444404409848
$ws.Range("T5").Value = 0.2189444404409849
$ws.Range("I6").Value = 0.2440410104700376
$ws.Range("J6").Value = 0.2440410104700377
$ws.Range("M6").Value = 2.309992333333333
$ws.Range("N6").Value = 6.929977
$ws.Range("O6").Value = 0.03775393480250816
$ws.Range("P6").Value = 0.03775393480250816
$ws.Range("S6").Value = 0.009213508398424011
$ws.Range("T6").Value = 0.009213508398424013
$ws.Range("I7").Value = 0.2440410104700376
$ws.Range("J7").Value = 0.2440410104700377
$ws.Range("M7").Value = 3.982169333333333
$ws.Range("N7").Value = 11.946508
$ws.Range("O7").Value = 0.06508357591224938
$ws.Range("P7").Value = 0.06508357591224936
$ws.Range("Q7").Value = 95.20558230147377
$ws.Range("R7").Value = 856.8502407132639
$ws.Range("S7").Value = 0.01588306163062874
$ws.Range("T7").Value = 0.01588306163062874
$ws.Range("G8").Value = 47.93131266666666
$ws.Range("H8").Value = 143.793938
$ws.Range("I8").Value = 0.489259703040525
$ws.Range("J8").Value = 0.4892597030405251
$ws.Range("M8").Value = 54.89331066666666
$ws.Range("N8").Value = 164.679932
$ws.Range("O8").Value = 0.8971624892852424
$ws.Range("P8").Value = 0.8971624892852424
$ws.Range("Q8").Value = 2631.108436872468
$ws.Range("R8").Value = 23679.97593185221
$ws.Range("S8").Value = 0.4389454530867959
$ws.Range("T8").Value = 0.4389454530867961
$ws.Range("G9").Value = 47.93131266666666
$ws.Range("H9").Value = 143.793938
$ws.Range("I9").Value = 0.489259703040525
$ws.Range("J9").Value = 0.4892597030405251
$ws.Range("M9").Value = 2.309992333333333
$ws.Range("N9").Value = 6.929977
$ws.Range("O9").Value = 0.03775393480250816
$ws.Range("P9").Value = 0.03775393480250816
$ws.Range("Q9").Value = 110.7209647866029
$ws.Range("R9").Value = 996.488683079426
$ws.Range("S9").Value = 0.01847147893008649
$ws.Range("T9").Value = 0.01847147893008649
$ws.Range("G10").Value = 47.93131266666666
$ws.Range("H10").Value = 143.793938
$ws.Range("I10").Value = 0.489259703040525
$ws.Range("J10").Value = 0.4892597030405251
$ws.Range("M10").Value = 3.982169333333333
$ws.Range("N10").Value = 11.946508
$ws.Range("O10").Value = 0.06508357591224938
$ws.Range("P10").Value = 0.06508357591224936
$ws.Range("Q10").Value = 190.8706034076115
$ws.Range("R10").Value = 1717.835430668504
$ws.Range("S10").Value = 0.0318427710236426
$ws.Range("T10").Value = 0.0318427710236426
